# Updated as of 3/22
# Refresh task names / dates / statuses on the "Basic Manual Gantt Chart"
# sheet to reflect the latest project schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task name (column B) updates ---------------------------------------

# Row 33: the "Dan & Tim" online cafe talk became just "Dan"
$ws.Range("B33").Value = "Online Café Talk (Dan)"

# Row 34: the "Sudip" online cafe talk gained "& Tim"
$ws.Range("B34").Value = "Online Café Talk (Sudip & Tim)"

# Row 35: "Meeting to Introduce John" was replaced by "Meeting # 13"
$ws.Range("B35").Value = "Meeting # 13"

# Row 37: "Meeting 13" renumbered to "Meeting #14"
$ws.Range("B37").Value = "Meeting #14"

# Row 36: "Meeting with Advisor # 5" now calls out Pablo by name
$ws.Range("B36").Value = "Meeting with Advisor (Pablo) # 5"

# --- Row 36 schedule moved earlier (3/27/2018 -> 3/24/2018) -------------

$ws.Range("C36").Value = 43183
$ws.Range("D36").Value = 43183

# --- Status updates (column F): mark rows 33-35 as Completed ------------
# Copy the formatting used by the already-"Completed" rows (e.g. F8) so the
# green fill/border that Excel applies for a completed task carries over,
# then set the text to "Completed".

$ws.Range("F8").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("F35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F33").Value = "Completed"
$ws.Range("F34").Value = "Completed"
$ws.Range("F35").Value = "Completed"
